# Applies odds updates to Sheet1 of the Betfair Back/Lay workbook for 2026-02-18.
# Each line below mirrors one <c>/<v> change in the source diff, updating the
# corresponding cell on the active worksheet to its new numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.4
$ws.Range("T2").Value = 1.72
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 980
$ws.Range("AB2").Value = 13.5
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 13.5
$ws.Range("AE2").Value = 980
$ws.Range("AF2").Value = 980
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 980
$ws.Range("AK2").Value = 980
$ws.Range("AN2").Value = 980
$ws.Range("AO2").Value = 980
$ws.Range("F3").Value = 2.52
$ws.Range("H3").Value = 2.62
$ws.Range("I3").Value = 2.84
$ws.Range("K3").Value = 4.3
$ws.Range("N3").Value = 2.76
$ws.Range("P3").Value = 1.85
$ws.Range("Q3").Value = 1.73
$ws.Range("Z3").Value = 980
$ws.Range("AD3").Value = 15
$ws.Range("AH3").Value = 1000
$ws.Range("F4").Value = 2.36
$ws.Range("G4").Value = 2.7
$ws.Range("H4").Value = 2.74
$ws.Range("I4").Value = 3.1
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 4.2
$ws.Range("P4").Value = 2.06
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 1.43
$ws.Range("S4").Value = 2.86
$ws.Range("U4").Value = 2.26
$ws.Range("V4").Value = 1.48
$ws.Range("W4").Value = 1.58
$ws.Range("X4").Value = 980
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 980
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 980
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 980
$ws.Range("G5").Value = 8.6
$ws.Range("P5").Value = 2.32
$ws.Range("U5").Value = 2.02
$ws.Range("AN5").Value = 140
$ws.Range("G6").Value = 5.1
$ws.Range("H6").Value = 1.74
$ws.Range("J6").Value = 4.1
$ws.Range("N6").Value = 5.7
$ws.Range("O6").Value = 1.18
$ws.Range("P6").Value = 2.58
$ws.Range("Q6").Value = 1.55
$ws.Range("R6").Value = 1.64
$ws.Range("S6").Value = 2.34
$ws.Range("W6").Value = 1.25
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 980
$ws.Range("AI6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 980
$ws.Range("AN6").Value = 980
$ws.Range("I7").Value = 1.88
$ws.Range("Q7").Value = 1.77
$ws.Range("R7").Value = 1.5
$ws.Range("W8").Value = 2.44
$ws.Range("Q9").Value = 2.06
$ws.Range("AB9").Value = 10
$ws.Range("AG9").Value = 11
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 2.28
$ws.Range("S10").Value = 2.94
$ws.Range("X11").Value = 14
$ws.Range("AB12").Value = 19
$ws.Range("I13").Value = 1.81
$ws.Range("T13").Value = 1.61